$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the Relatorio sheet (sheet2) content: replace the old "report"
#    table with the same Produto/Quantidade/Preco Unitario/Total layout used
#    on Vendas, but with a new set of rows (no totals filled in column D).
# ---------------------------------------------------------------------------
$wsRelatorio = $wb.Worksheets.Item("Relatorio")

$wsRelatorio.Range("A1").Value = "Produto"
$wsRelatorio.Range("B1").Value = "Quantidade"
$wsRelatorio.Range("C1").Value = "Preço Unitário"
$wsRelatorio.Range("D1").Value = "Total"

$wsRelatorio.Range("A2").Value = "Kiwi"
$wsRelatorio.Range("B2").Value = 16
$wsRelatorio.Range("C2").Value = 7.3

$wsRelatorio.Range("A3").Value = "Morango"
$wsRelatorio.Range("B3").Value = 13
$wsRelatorio.Range("C3").Value = 7.14

$wsRelatorio.Range("A4").Value = "Uva"
$wsRelatorio.Range("B4").Value = 20
$wsRelatorio.Range("C4").Value = 8.38

$wsRelatorio.Range("A5").Value = "Maçã"
$wsRelatorio.Range("B5").Value = 7
$wsRelatorio.Range("C5").Value = 2.61

$wsRelatorio.Range("A6").Value = "Manga"
$wsRelatorio.Range("B6").Value = 8
$wsRelatorio.Range("C6").Value = 7.46

# ---------------------------------------------------------------------------
# 2. Add the three monthly sheets + the summary sheet, in order, at the end
#    of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsJaneiro = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsJaneiro.Name = "Janeiro"

$wsFevereiro = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsJaneiro)
$wsFevereiro.Name = "Fevereiro"

$wsMarco = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsFevereiro)
$wsMarco.Name = "Março"

$wsResumo = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsMarco)
$wsResumo.Name = "Resumo"

# --- Janeiro ---
$wsJaneiro.Range("A1").Value = "Produto"
$wsJaneiro.Range("B1").Value = "Quantidade"
$wsJaneiro.Range("C1").Value = "Preço Unitário"
$wsJaneiro.Range("D1").Value = "Total"

$wsJaneiro.Range("A2").Value = "Uva"
$wsJaneiro.Range("B2").Value = 11
$wsJaneiro.Range("C2").Value = 5.14

$wsJaneiro.Range("A3").Value = "Kiwi"
$wsJaneiro.Range("B3").Value = 3
$wsJaneiro.Range("C3").Value = 1.43

$wsJaneiro.Range("A4").Value = "Manga"
$wsJaneiro.Range("B4").Value = 18
$wsJaneiro.Range("C4").Value = 5.19

$wsJaneiro.Range("A5").Value = "Morango"
$wsJaneiro.Range("B5").Value = 7
$wsJaneiro.Range("C5").Value = 6.95

$wsJaneiro.Range("A6").Value = "Maçã"
$wsJaneiro.Range("B6").Value = 1
$wsJaneiro.Range("C6").Value = 3.34

# --- Fevereiro ---
$wsFevereiro.Range("A1").Value = "Produto"
$wsFevereiro.Range("B1").Value = "Quantidade"
$wsFevereiro.Range("C1").Value = "Preço Unitário"
$wsFevereiro.Range("D1").Value = "Total"

$wsFevereiro.Range("A2").Value = "Maçã"
$wsFevereiro.Range("B2").Value = 19
$wsFevereiro.Range("C2").Value = 1.73

$wsFevereiro.Range("A3").Value = "Manga"
$wsFevereiro.Range("B3").Value = 20
$wsFevereiro.Range("C3").Value = 9.68

$wsFevereiro.Range("A4").Value = "Morango"
$wsFevereiro.Range("B4").Value = 19
$wsFevereiro.Range("C4").Value = 9.83

$wsFevereiro.Range("A5").Value = "Uva"
$wsFevereiro.Range("B5").Value = 11
$wsFevereiro.Range("C5").Value = 4.37

$wsFevereiro.Range("A6").Value = "Kiwi"
$wsFevereiro.Range("B6").Value = 14
$wsFevereiro.Range("C6").Value = 6.88

# --- Março ---
$wsMarco.Range("A1").Value = "Produto"
$wsMarco.Range("B1").Value = "Quantidade"
$wsMarco.Range("C1").Value = "Preço Unitário"
$wsMarco.Range("D1").Value = "Total"

$wsMarco.Range("A2").Value = "Manga"
$wsMarco.Range("B2").Value = 13
$wsMarco.Range("C2").Value = 5.41

$wsMarco.Range("A3").Value = "Morango"
$wsMarco.Range("B3").Value = 20
$wsMarco.Range("C3").Value = 8.63

$wsMarco.Range("A4").Value = "Uva"
$wsMarco.Range("B4").Value = 3
$wsMarco.Range("C4").Value = 1.25

$wsMarco.Range("A5").Value = "Maçã"
$wsMarco.Range("B5").Value = 8
$wsMarco.Range("C5").Value = 7.05

$wsMarco.Range("A6").Value = "Kiwi"
$wsMarco.Range("B6").Value = 20
$wsMarco.Range("C6").Value = 5.49

# --- Resumo (headers only) ---
$wsResumo.Range("A1").Value = "Produto"
$wsResumo.Range("B1").Value = "Quantidade"
$wsResumo.Range("C1").Value = "Preço Unitário"
$wsResumo.Range("D1").Value = "Total"

# ---------------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping to match the target state:
#    Relatorio becomes the active sheet (activeTab=1); Vendas keeps a
#    lingering selection at P19 but is no longer the active tab.
# ---------------------------------------------------------------------------
$wsVendas = $wb.Worksheets.Item("Vendas")
$wsVendas.Activate()
$wsVendas.Range("P19").Select()

$wsRelatorio.Activate()
$wsRelatorio.Range("A1").Select()
